$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric cells -----------------------------------------------
$ws.Range("A3").Value = 131082382
$ws.Range("B3").Value = 58043
$ws.Range("E3").Value = 103021
$ws.Range("Q3").Value = 572824
$ws.Range("R3").Value = 6447033
$ws.Range("S3").Value = 10

# --- Plain text cells ----------------------------------------------------
$ws.Range("D3").Value = "NT"
$ws.Range("F3").Value = "Talltita"
$ws.Range("G3").Value = "Poecile montanus"
$ws.Range("H3").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("M3").Value = "permanent revir"
$ws.Range("P3").Value = "Öjsjön-Falerum, Sm"
$ws.Range("T3").Value = "Östergötland"
$ws.Range("U3").Value = "Åtvidaberg"
$ws.Range("V3").Value = "Småland"
$ws.Range("W3").Value = "Gärdserum"
$ws.Range("AC3").Value = "Revirparet i delområde av sitt revir"
$ws.Range("AW3").Value = "Steve Daurer"
$ws.Range("AX3").Value = "Steve Daurer"

# --- Text cells whose content looks numeric/date-like --------------------
# A direct .Value assignment would be auto-coerced to a number by Excel,
# so the value is entered as a formula producing the literal text and
# then converted in place to a plain (non-formula) text value.
$ws.Range("I3").Formula = '="2"'
$ws.Range("I3").Copy()
$ws.Range("I3").PasteSpecial(-4163)

$ws.Range("Y3").Formula = '="2026-01-31"'
$ws.Range("Y3").Copy()
$ws.Range("Y3").PasteSpecial(-4163)

$ws.Range("AA3").Formula = '="2026-01-31"'
$ws.Range("AA3").Copy()
$ws.Range("AA3").PasteSpecial(-4163)

# --- Boolean cells ---------------------------------------------------------
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
